# "Set rendimientos 13 May" - update the raw Resueltos/Pendientes figures
# on the "tareas" sheet. D (Tareas totales), E (Media de Trabajos), F
# (Rendimiento) and G (Porcentaje) are all formulas, so they recompute
# automatically once B/C change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tareas")

$ws.Range("B2").Value = 116

$ws.Range("B3").Value = 29

$ws.Range("B4").Value = 41
$ws.Range("C4").Value = 5

$ws.Range("B6").Value = 14
$ws.Range("C6").Value = 1

$ws.Range("C7").Value = 6

$ws.Range("B8").Value = 48
$ws.Range("C8").Value = 16

$ws.Range("B9").Value = 21
$ws.Range("C9").Value = 2

$ws.Range("C11").Value = 10

$ws.Range("B12").Value = 26
$ws.Range("C12").Value = 0

$ws.Range("C13").Value = 4

$ws.Range("B16").Value = 35

$ws.Range("B17").Value = 21
$ws.Range("C17").Value = 5

$ws.Range("B18").Value = 9

$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 1

$ws.Range("B25").Value = 6

$ws.Range("B27").Value = 8
$ws.Range("C27").Value = 5

$ws.Range("B28").Value = 8
$ws.Range("C28").Value = 1

# Match the saved view state: scrolled down with C29 selected, zoomed to 130%.
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 130
